# Update automatico via Actualizar 02-06-2021 03-07-22
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Refresh the "last checked" timestamp for the previous batch of checks
#    (rows 954-967) to the latest re-check time.
# ---------------------------------------------------------------------------
$updatedTimestamp = 44233.10887606481
for ($r = 954; $r -le 967; $r++) {
    $ws.Cells.Item($r, 4).Value = $updatedTimestamp
}

# ---------------------------------------------------------------------------
# 2. Append a brand-new batch of availability checks (rows 968-981), one row
#    per monitored service, all sharing the same check timestamp.
# ---------------------------------------------------------------------------
$newTimestamp = 44233.13005763067

# Each entry: row, service name, display URL text, hyperlink target,
#             hyperlink sub-address (fragment), row number in A column style
$entries = @(
    @{ Row = 968; Name = "Odoo";               Url = "https://www.dataintelligence-group.com/";                          Target = "https://www.dataintelligence-group.com/";                          SubAddress = "" },
    @{ Row = 969; Name = "Blackbox";            Url = "https://serviciodashboard.azurewebsites.net/";                     Target = "https://serviciodashboard.azurewebsites.net/";                     SubAddress = "" },
    @{ Row = 970; Name = "PowerBI";             Url = "https://powerbi.microsoft.com/es-es/";                            Target = "https://powerbi.microsoft.com/es-es/";                             SubAddress = "" },
    @{ Row = 971; Name = "Dropbox";             Url = "https://www.dropbox.com/";                                        Target = "https://www.dropbox.com/";                                         SubAddress = "" },
    @{ Row = 972; Name = "Odoo";                Url = "https://dataintelligence.store/";                                 Target = "https://dataintelligence.store/";                                  SubAddress = "" },
    @{ Row = 973; Name = "GEE";                 Url = "https://app-data-i.users.earthengine.app/";                       Target = "https://app-data-i.users.earthengine.app/";                        SubAddress = "" },
    @{ Row = 974; Name = "UtilidadesOdoo";      Url = "https://odooutil.azurewebsites.net/";                             Target = "https://odooutil.azurewebsites.net/";                              SubAddress = "" },
    @{ Row = 975; Name = "Filtros Dashboard";   Url = "https://filtradordashboard.azurewebsites.net/";                   Target = "https://filtradordashboard.azurewebsites.net/";                    SubAddress = "" },
    @{ Row = 976; Name = "MapStore";            Url = "https://ide.dataintelligence-group.com/mapstore/#/";              Target = "https://ide.dataintelligence-group.com/mapstore/";                 SubAddress = "/" },
    @{ Row = 977; Name = "GeoServer";           Url = "https://ide.dataintelligence-group.com/geoserver/web/?0";         Target = "https://ide.dataintelligence-group.com/geoserver/web/?0";          SubAddress = "" },
    @{ Row = 978; Name = "Tomcat";              Url = "https://ide.dataintelligence-group.com/";                         Target = "https://ide.dataintelligence-group.com/";                          SubAddress = "" },
    @{ Row = 979; Name = "Shiny";               Url = "https://rpubs.com/dataintelligence/";                             Target = "https://rpubs.com/dataintelligence/";                              SubAddress = "" },
    @{ Row = 980; Name = "Github";              Url = "https://github.com/Sud-Austral/";                                 Target = "https://github.com/Sud-Austral/";                                  SubAddress = "" },
    @{ Row = 981; Name = "EZ Exporter";         Url = "https://ezexporter.highviewapps.com/exports/export-profile/";     Target = "https://ezexporter.highviewapps.com/exports/export-profile/";     SubAddress = "" }
)

foreach ($entry in $entries) {
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value = $entry.Name
    $ws.Cells.Item($r, 2).Value = $entry.Url
    $ws.Cells.Item($r, 3).Value = "Disponible"
    $ws.Cells.Item($r, 4).Value = $newTimestamp
    # Keep the same date/time display format used throughout column D.
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    if ($entry.SubAddress -ne "") {
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $entry.Target, $entry.SubAddress)
    } else {
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $entry.Target)
    }
    # Restore the standard hyperlink look (matches the styling used by every
    # other URL cell in column B) after Hyperlinks.Add applies its own style.
    $ws.Cells.Item($r, 2).Style = "Hyperlink"
}
